$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.829
$ws.Range("E4").Value = 13.195
$ws.Range("E5").Value = 13.654
$ws.Range("C6").Value = -12.715
$ws.Range("E6").Value = 13.493
$ws.Range("C7").Value = -13.214
$ws.Range("C8").Value = -12.694
$ws.Range("E8").Value = 13.423
$ws.Range("C16").Value = -11.738
$ws.Range("E16").Value = 13.094
$ws.Range("C20").Value = -13.349
$ws.Range("C21").Value = -13.214
$ws.Range("E22").Value = 13.238
